$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row74
$ws.Range("H74").Value = 3891.0908
$ws.Range("I74").Value = 4006.9333
$ws.Range("J74").Value = 3642.8572
$ws.Range("K74").Value = 4006.9333
$ws.Range("L74").Value = 3642.8572
$ws.Range("M74").Value = -3070.9333
$ws.Range("N74").Value = -5514.8572

# ALC!row77
$ws.Range("H77").Value = 3891.0908
$ws.Range("I77").Value = 4006.9333
$ws.Range("J77").Value = 3642.8572
$ws.Range("K77").Value = 20034.6665
$ws.Range("L77").Value = 18214.286
$ws.Range("M77").Value = -15354.6665
$ws.Range("N77").Value = -27574.286

# ALC!row132
$ws.Range("H132").Value = 847268.25
$ws.Range("I132").Value = 2696.2703
$ws.Range("J132").Value = 2335323.8
$ws.Range("K132").Value = 8088.8109
$ws.Range("L132").Value = 7005971.399999999
$ws.Range("M132").Value = -5558.8109
$ws.Range("N132").Value = -7011031.399999999

# ALC!row135
$ws.Range("H135").Value = 19611.928
$ws.Range("I135").Value = 23277.955
$ws.Range("J135").Value = 4614.5454
$ws.Range("K135").Value = 209501.595
$ws.Range("L135").Value = 41530.9086
$ws.Range("M135").Value = -206966.595
$ws.Range("N135").Value = -46600.9086

# ALC!row138
$ws.Range("H138").Value = 2316797.8
$ws.Range("I138").Value = 1399.4889
$ws.Range("J138").Value = 6175795
$ws.Range("K138").Value = 4198.4667
$ws.Range("L138").Value = 18527385
$ws.Range("M138").Value = 941.5333000000001
$ws.Range("N138").Value = -18537665

$ws = $wb.Worksheets.Item("ARM")
# ARM!row32
$ws.Range("H32").Value = 17785
$ws.Range("I32").Value = 17104.574
$ws.Range("J32").Value = 20470.895
$ws.Range("K32").Value = 17104.574
$ws.Range("L32").Value = 20470.895
$ws.Range("M32").Value = -16817.574
$ws.Range("N32").Value = -21044.895

# ARM!row33
$ws.Range("H33").Value = 5864.5
$ws.Range("I33").Value = 2700
$ws.Range("J33").Value = 9029
$ws.Range("K33").Value = 2700
$ws.Range("L33").Value = 9029
$ws.Range("M33").Value = -2371
$ws.Range("N33").Value = -9687

# ARM!row61
$ws.Range("H61").Value = 16162556
$ws.Range("I61").Value = 18201000
$ws.Range("J61").Value = 146202
$ws.Range("K61").Value = 18201000
$ws.Range("L61").Value = 146202
$ws.Range("M61").Value = -18200788
$ws.Range("N61").Value = -146626

# ARM!row74
$ws.Range("H74").Value = 10163208
$ws.Range("I74").Value = 12860239
$ws.Range("J74").Value = 145664.14
$ws.Range("K74").Value = 12860239
$ws.Range("L74").Value = 145664.14
$ws.Range("M74").Value = -12859365
$ws.Range("N74").Value = -147412.14

# ARM!row77
$ws.Range("H77").Value = 10163208
$ws.Range("I77").Value = 12860239
$ws.Range("J77").Value = 145664.14
$ws.Range("K77").Value = 64301195
$ws.Range("L77").Value = 728320.7000000001
$ws.Range("M77").Value = -64296827
$ws.Range("N77").Value = -737056.7000000001

# ARM!row136
$ws.Range("H136").Value = 16162556
$ws.Range("I136").Value = 18201000
$ws.Range("J136").Value = 146202
$ws.Range("K136").Value = 54603000
$ws.Range("L136").Value = 438606
$ws.Range("M136").Value = -54600450
$ws.Range("N136").Value = -443706

$ws = $wb.Worksheets.Item("BSM")
# BSM!row7
$ws.Range("H7").Value = 645.75
$ws.Range("I7").Value = 301.5
$ws.Range("J7").Value = 990
$ws.Range("K7").Value = 301.5
$ws.Range("L7").Value = 990
$ws.Range("M7").Value = -188.5
$ws.Range("N7").Value = -1216

# BSM!row107
$ws.Range("H107").Value = 3615.8572
$ws.Range("I107").Value = 3768.5
$ws.Range("J107").Value = 2700
$ws.Range("K107").Value = 3768.5
$ws.Range("L107").Value = 2700
$ws.Range("M107").Value = -1848.5
$ws.Range("N107").Value = -6540

# BSM!row134
$ws.Range("H134").Value = 1699.4906
$ws.Range("I134").Value = 1094.4103
$ws.Range("J134").Value = 3385.0715
$ws.Range("K134").Value = 3283.2309
$ws.Range("L134").Value = 10155.2145
$ws.Range("M134").Value = -748.2309
$ws.Range("N134").Value = -15225.2145

$ws = $wb.Worksheets.Item("CRP")
# CRP!row32
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

# CRP!row132
$ws.Range("H132").Value = 18983.896
$ws.Range("I132").Value = 1565.1951
$ws.Range("J132").Value = 60993.707
$ws.Range("K132").Value = 4695.5853
$ws.Range("L132").Value = 182981.121
$ws.Range("M132").Value = -2165.5853
$ws.Range("N132").Value = -188041.121

# CRP!row134
$ws.Range("H134").Value = 15364.474
$ws.Range("I134").Value = 1104.9482
$ws.Range("K134").Value = 3314.8446
$ws.Range("M134").Value = -779.8446000000004

$ws = $wb.Worksheets.Item("CUL")
# CUL!row23
$ws.Range("H23").Value = 120.210526
$ws.Range("I23").Value = 101.833336
$ws.Range("J23").Value = 128.6923
$ws.Range("K23").Value = 305.500008
$ws.Range("L23").Value = 386.0769
$ws.Range("M23").Value = -70.50000799999998
$ws.Range("N23").Value = -856.0769

# CUL!row129
$ws.Range("H129").Value = 2978081.2
$ws.Range("I129").Value = 1366.8462
$ws.Range("J129").Value = 5557900.5
$ws.Range("K129").Value = 4100.5386
$ws.Range("L129").Value = 16673701.5
$ws.Range("M129").Value = 899.4614000000001
$ws.Range("N129").Value = -16683701.5

# CUL!row131
$ws.Range("H131").Value = 1122.8628
$ws.Range("J131").Value = 1232.6364
$ws.Range("L131").Value = 3697.9092
$ws.Range("N131").Value = -13777.9092

$ws = $wb.Worksheets.Item("LTW")
# LTW!row40
$ws.Range("H40").Value = 3226.7666
$ws.Range("I40").Value = 2774.0435
$ws.Range("K40").Value = 2774.0435
$ws.Range("M40").Value = -2638.0435

# LTW!row132
$ws.Range("H132").Value = 23022.426
$ws.Range("I132").Value = 1403.1428
$ws.Range("K132").Value = 4209.428400000001
$ws.Range("M132").Value = -1679.428400000001

# LTW!row136
$ws.Range("H136").Value = 1259.1111
$ws.Range("I136").Value = 884.5714
$ws.Range("J136").Value = 2570
$ws.Range("K136").Value = 2653.7142
$ws.Range("L136").Value = 7710
$ws.Range("M136").Value = -103.7142000000003
$ws.Range("N136").Value = -12810

$ws = $wb.Worksheets.Item("WVR")
# WVR!row62
$ws.Range("H62").Value = 2900.55
$ws.Range("I62").Value = 2906.375
$ws.Range("J62").Value = 2877.25
$ws.Range("K62").Value = 2906.375
$ws.Range("L62").Value = 2877.25
$ws.Range("M62").Value = -2282.375
$ws.Range("N62").Value = -4125.25

# WVR!row65
$ws.Range("H65").Value = 2900.55
$ws.Range("I65").Value = 2906.375
$ws.Range("J65").Value = 2877.25
$ws.Range("K65").Value = 14531.875
$ws.Range("L65").Value = 14386.25
$ws.Range("M65").Value = -11411.875
$ws.Range("N65").Value = -20626.25

# WVR!row122
$ws.Range("H122").Value = 1862.6586
$ws.Range("J122").Value = 2687.4375
$ws.Range("L122").Value = 8062.3125
$ws.Range("N122").Value = -12962.3125

# WVR!row132
$ws.Range("H132").Value = 47738.723
$ws.Range("I132").Value = 34200.066
$ws.Range("J132").Value = 78981.766
$ws.Range("K132").Value = 102600.198
$ws.Range("L132").Value = 236945.298
$ws.Range("M132").Value = -100070.198
$ws.Range("N132").Value = -242005.298

# WVR!row136
$ws.Range("H136").Value = 55279.42
$ws.Range("I136").Value = 46710.816
$ws.Range("J136").Value = 67061.25
$ws.Range("K136").Value = 140132.448
$ws.Range("L136").Value = 201183.75
$ws.Range("M136").Value = -137582.448
$ws.Range("N136").Value = -206283.75
